$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-4) need to be re-sorted chronologically by date (column D).
# Resulting values per row for columns D, J, K, L, M, P:
$ws.Range("D2").Value = 44175
$ws.Range("J2").Value = 1400
$ws.Range("K2").Value = 1900
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = 1950
$ws.Range("P2").Value = 1950

$ws.Range("D3").Value = 44200
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 1400
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = 1450
$ws.Range("P3").Value = 1450

$ws.Range("D4").Value = 44210
$ws.Range("J4").Value = 1450
$ws.Range("K4").Value = 1600
$ws.Range("L4").Value = 1700
$ws.Range("M4").Value = 1650
$ws.Range("P4").Value = 1650

$wb.Save()
